# daily auto push: 2025-10-02 13:35 UTC
# Append a new data row (row 52) to the bottom of the log table on the
# active sheet, mirroring the existing rows (A: date text, B: weekday
# text, C/D: numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 52

# Column A holds a date-like string ("2025/10/02"). Assigning such a
# string directly to .Value makes Excel auto-convert it into a date
# serial number, so we force the cell to a text format first and then
# restore the cell style to "Normal" afterwards so no stray style index
# is left on the cell (matching the unstyled data rows above it).
$ws.Range("A$newRow").NumberFormat = "@"
$ws.Range("A$newRow").Value = "2025/10/02"
$ws.Range("A$newRow").Style = "Normal"

$ws.Range("B$newRow").Value = "木"
$ws.Range("C$newRow").Value = 20
$ws.Range("D$newRow").Value = 27
